$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: attendance date headers H4/I4 switch from date serials to text dates ---
$ws.Range("H4").Value = "17/9/2022"
$ws.Range("I4").Value = "24/9/2022"

# --- Clear the S.NO column (A5:A41) roll-number sequence values ---
$snoRows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,26,27,28,29,30,31,32,34,35,37,39,40,41)
foreach ($r in $snoRows) {
    $ws.Cells.Item($r, 1).ClearContents()
}

# --- Populate the new attendance column I (column 9) for the 24/9/2022 session ---
$iValues = @{
    5  = 3;  6  = 3;  7  = 0;  8  = 3;  9  = 3;
    10 = 3;  11 = 3;  12 = 3;  13 = 0;  14 = 0;
    15 = 3;  16 = 0;  17 = 3;  18 = 3;  19 = 3;
    20 = 0;  21 = 3;  22 = 0;  23 = 0;  24 = 0;
    25 = 0;  26 = 0;  27 = 0;  28 = 0;  29 = 0;
    30 = 0;  31 = 0;  32 = 3;  33 = 0;  34 = 3;
    35 = 3;  36 = 3;  37 = 3;  38 = 3;  39 = 3;
    40 = 0;  41 = 0
}
foreach ($r in 5..41) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
}

# --- Extend the weekly-average formula row to include the new column I ---
$ws.Range("I42").Formula = "=SUM(I5:I41)/3"

# --- Update the active selection to reflect where the editor left the cursor ---
$ws.Range("F5").Select()
